$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294; this shifts the existing rows 294-397 down to 295-398
$ws.Rows(294).Insert()

# Populate the newly inserted row 294 with the new data record
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44900
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100114014
$ws.Range("G294").Value = "Betarraga"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 250
$ws.Range("K294").Value = 1200
$ws.Range("L294").Value = 1200
$ws.Range("M294").Value = 1200
$ws.Range("N294").Value = '$/paquete 5 unidades'
$ws.Range("O294").Value = "Región del Maule"
$ws.Range("P294").Value = 240
$ws.Range("Q294").Value = 5
$ws.Range("R294").Value = "Hortaliza"
